$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exposure conditions")

# Delete the last 4 rows (rows 8-11), which removes the duplicate
# "CONTROL (SEE VEHICLE)" replicate 3/4 rows and the extra pair of
# "EXTRACTION BLANK" rows, shrinking the data range from A1:L11 to A1:L7.
$ws.Rows("8:11").Delete()

# Rows 2-5 previously described "chemical1" doses at TP1; they now
# represent the control (see vehicle) rows at TP0.
$ws.Range("J2").Value = "CONTROL (SEE VEHICLE)"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = "TP0"

$ws.Range("J3").Value = "CONTROL (SEE VEHICLE)"
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "TP0"

$ws.Range("J4").Value = "CONTROL (SEE VEHICLE)"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "TP0"

$ws.Range("J5").Value = "CONTROL (SEE VEHICLE)"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = "TP0"

# Rows 6-7 (previously CONTROL (SEE VEHICLE) replicate 1/2) now become
# the EXTRACTION BLANK rows at TP0, matching the rows that used to be
# rows 10-11 before the deletion above. Their dose column ("0") is
# stored as text, so force the cell format to text before assigning.
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "EXTRACTION BLANK"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "0"
$ws.Range("L6").Value = "TP0"

$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "EXTRACTION BLANK"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "0"
$ws.Range("L7").Value = "TP0"
